$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptocurrency market data (price & volume/1h changes,
# plus a few re-ordered rows) while preserving original text formatting.

$ws.Cells.Item(2, 4).NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = '66.988.03'
$ws.Cells.Item(2, 5).Value = '  +0.60%  '
$ws.Cells.Item(3, 4).NumberFormat = "@"
$ws.Cells.Item(3, 4).Value = '3.806.01'
$ws.Cells.Item(3, 5).Value = '  -1.08%  '
$ws.Cells.Item(4, 4).NumberFormat = "@"
$ws.Cells.Item(4, 4).Value = '0.998'
$ws.Cells.Item(4, 5).Value = '  -0.20%  '
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '436.20'
$ws.Cells.Item(5, 5).Value = '  +1.47%  '
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '139.57'
$ws.Cells.Item(6, 5).Value = '  +6.50%  '
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = '0.623'
$ws.Cells.Item(7, 5).Value = '  +2.05%  '
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = '0.998'
$ws.Cells.Item(8, 5).Value = '  -0.05%  '
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = '0.736'
$ws.Cells.Item(9, 5).Value = '  +1.21%  '
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = '0.155'
$ws.Cells.Item(10, 5).Value = '  -7.22%  '
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = '0.0000320'
$ws.Cells.Item(11, 5).Value = '  -11.87%  '
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = '43.03'
$ws.Cells.Item(12, 5).Value = '  +5.50%  '
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = '10.42'
$ws.Cells.Item(13, 5).Value = '  +3.65%  '
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = '4.396.95'
$ws.Cells.Item(14, 5).Value = '  -1.03%  '
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = '14.90'
$ws.Cells.Item(15, 5).Value = '  -5.44%  '
$ws.Cells.Item(16, 5).Value = '  -0.42%  '
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = '3.788.52'
$ws.Cells.Item(17, 5).Value = '  -1.81%  '
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = '19.96'
$ws.Cells.Item(18, 5).Value = '  +1.93%  '
$ws.Cells.Item(19, 5).Value = '  +6.86%  '
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = '67.017.61'
$ws.Cells.Item(20, 5).Value = '  +0.09%  '
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = '419.86'
$ws.Cells.Item(21, 5).Value = '  +2.67%  '
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '14.72'
$ws.Cells.Item(22, 5).Value = '  +1.93%  '
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = '3.25'
$ws.Cells.Item(23, 5).Value = '  +6.83%  '
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = '86.26'
$ws.Cells.Item(24, 5).Value = '  +1.13%  '
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = '37.44'
$ws.Cells.Item(25, 5).Value = '  +1.34%  '
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = '3.40'
$ws.Cells.Item(26, 5).Value = '  +4.34%  '
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = '9.73'
$ws.Cells.Item(27, 5).Value = '  +35.85%  '
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = '5.58'
$ws.Cells.Item(28, 5).Value = '  -1.61%  '
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = '9.85'
$ws.Cells.Item(29, 5).Value = '  +3.76%  '
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = '738.39'
$ws.Cells.Item(30, 5).Value = '  +7.16%  '
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = '13.84'
$ws.Cells.Item(31, 5).Value = '  +11.10%  '
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = '0.134'
$ws.Cells.Item(32, 5).Value = '  +10.62%  '
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = '2.75'
$ws.Cells.Item(33, 5).Value = '  +2.71%  '
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = '43.29'
$ws.Cells.Item(34, 5).Value = '  +11.88%  '
$ws.Cells.Item(35, 2).Value = 'Kaspa'
$ws.Cells.Item(35, 3).Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = '0.159'
$ws.Cells.Item(35, 5).Value = '  +4.56%  '
$ws.Cells.Item(36, 2).Value = 'Dai'
$ws.Cells.Item(36, 3).Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = '1.00'
$ws.Cells.Item(36, 5).Value = '  +0.02%  '
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = '5.61'
$ws.Cells.Item(37, 5).Value = '  +23.90%  '
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = '56.36'
$ws.Cells.Item(38, 5).Value = '  +2.07%  '
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = '0.0482'
$ws.Cells.Item(39, 5).Value = '  +5.39%  '
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '2.73'
$ws.Cells.Item(40, 5).Value = '  +40.79%  '
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = '2.97'
$ws.Cells.Item(41, 5).Value = '  -3.77%  '
$ws.Cells.Item(42, 2).Value = 'PEPE'
$ws.Cells.Item(42, 3).Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = '0.0₃0687'
$ws.Cells.Item(42, 5).Value = '  -12.89%  '
$ws.Cells.Item(43, 2).Value = 'Stellar'
$ws.Cells.Item(43, 3).Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = '0.142'
$ws.Cells.Item(43, 5).Value = '  +3.75%  '
$ws.Cells.Item(44, 2).Value = 'ApeXProtocol'
$ws.Cells.Item(44, 3).Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = '3.33'
$ws.Cells.Item(44, 5).Value = '  +6.62%  '
$ws.Cells.Item(45, 2).Value = 'FirstDigitalUSD'
$ws.Cells.Item(45, 3).Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = '1.00'
$ws.Cells.Item(45, 5).Value = '  -0.10%  '
$ws.Cells.Item(46, 2).Value = 'TheGraph'
$ws.Cells.Item(46, 3).Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = '0.331'
$ws.Cells.Item(46, 5).Value = '  +12.76%  '
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = '3.32'
$ws.Cells.Item(47, 5).Value = '  +0.83%  '
$ws.Cells.Item(48, 2).Value = 'ARBITRUM'
$ws.Cells.Item(48, 3).Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = '2.09'
$ws.Cells.Item(48, 5).Value = '  +0.61%  '
$ws.Cells.Item(49, 2).Value = 'WEMIXToken'
$ws.Cells.Item(49, 3).Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = '2.64'
$ws.Cells.Item(49, 5).Value = '  +4.31%  '
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = '142.70'
$ws.Cells.Item(50, 5).Value = '  -3.79%  '
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = '2.84'
$ws.Cells.Item(51, 5).Value = '  +1.92%  '
